$d = $word.ActiveDocument

# The original paragraph is made of three separately-formatted runs:
#   "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe ("
#   "http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/"   (styled as a hyperlink)
#   ")."
# The edit collapses them into a single, unformatted run with the updated
# (2022) year, so we locate the whole sentence, delete it, and retype the
# updated text as one plain run.
$old = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$new = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $rng.Text -eq $old) {
    $rng.Delete()
    $rng.InsertBefore($new)
}
